$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("plate")

$ws.Range("A2").Value = "CHODERA20200226a1.itc"
$ws.Range("A3").Value = "CHODERA20200226a2.itc"
$ws.Range("A4").Value = "CHODERA20200226a3.itc"
$ws.Range("A5").Value = "CHODERA20200226a4.itc"
$ws.Range("A6").Value = "CHODERA20200226a5.itc"
$ws.Range("A7").Value = "CHODERA20200226a6.itc"
$ws.Range("A8").Value = "CHODERA20200226a7.itc"
$ws.Range("A9").Value = "CHODERA20200226a8.itc"
$ws.Range("A10").Value = "CHODERA20200226a9.itc"
$ws.Range("A11").Value = "CHODERA20200226a10.itc"
